$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.193.39"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "3.595.82"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.43"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.19"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").Value = "3.584.90"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.219"
$ws.Range("E10").Value = "  +17.55%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.45"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000322"
$ws.Range("E13").Value = "  +6.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.55"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "4.171.14"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "71.292.87"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.33"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "3.594.18"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.38"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "567.48"
$ws.Range("E20").Value = "  +6.66%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.62"
$ws.Range("E23").Value = "  -10.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.12"
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.61"
$ws.Range("E25").Value = "  +4.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.03"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.36"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.41"
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.33"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.49"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.36"
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "554.90"
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.416"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").Value = "0.0₃0810"
$ws.Range("E38").Value = "  +4.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.75"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").Value = "3.515.39"
$ws.Range("E42").Value = "  +10.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.44"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0449"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.97"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.46"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.997"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  -1.73%  "
